# Generate Report for Handoff
#
# The localization status report has moved from "In Translation" to
# "Ready for handoff": update the Status cells on every sheet plus the
# two generation timestamps that accompany that state change, then
# re-size the columns that now hold the longer status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"              # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"              # de-de status
$wsOverview.Range("G2").Value = "2016-08-18 20:40:07"             # Latest HO Xliff Generate Date

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"                   # Status
$wsZhCn.Range("H2").Value = "2016-08-18 20:39:56"                 # Latest Handoff Datetime

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"                   # Status
$wsDeDe.Range("H2").Value = "2016-08-18 20:40:07"                 # Latest Handoff Datetime

# --- Resize the columns that now contain the longer status text ------
$wsOverview.Columns("E:F").ColumnWidth = 16.33
$wsZhCn.Columns("C:C").ColumnWidth = 16.33
$wsDeDe.Columns("C:C").ColumnWidth = 16.33
